$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.08367710560560226
$ws.Cells.Item(2, 2).Value = 0.9736980199813843
$ws.Cells.Item(2, 3).Value = 0.01179144158959389
$ws.Cells.Item(2, 4).Value = 0.9963904619216919
$ws.Cells.Item(3, 1).Value = 0.01313626766204834
$ws.Cells.Item(3, 2).Value = 0.9981400370597839
$ws.Cells.Item(3, 3).Value = 0.04806109145283699
$ws.Cells.Item(3, 4).Value = 0.9847771525382996
$ws.Cells.Item(4, 1).Value = 0.008491170592606068
$ws.Cells.Item(4, 2).Value = 0.9983546137809753
$ws.Cells.Item(4, 3).Value = 0.01407028362154961
$ws.Cells.Item(4, 4).Value = 0.9922316670417786
$ws.Cells.Item(5, 1).Value = 0.004161574877798557
$ws.Cells.Item(5, 2).Value = 0.9989507794380188
$ws.Cells.Item(5, 3).Value = 0.001539218705147505
$ws.Cells.Item(5, 4).Value = 0.9996076822280884
$ws.Cells.Item(6, 1).Value = 0.002075845841318369
$ws.Cells.Item(6, 2).Value = 0.9995469450950623
$ws.Cells.Item(6, 3).Value = 0.006949895527213812
$ws.Cells.Item(6, 4).Value = 0.9977244138717651
$ws.Cells.Item(7, 1).Value = 0.003393801860511303
$ws.Cells.Item(7, 2).Value = 0.9993323087692261
$ws.Cells.Item(7, 3).Value = 0.005113150924444199
$ws.Cells.Item(7, 4).Value = 0.9981167316436768
$ws.Cells.Item(8, 1).Value = 0.001545752165839076
$ws.Cells.Item(8, 2).Value = 0.9996423125267029
$ws.Cells.Item(8, 3).Value = 0.001036555273458362
$ws.Cells.Item(8, 4).Value = 0.999764621257782
$ws.Cells.Item(9, 1).Value = 0.001442248234525323
$ws.Cells.Item(9, 2).Value = 0.9997138381004333
$ws.Cells.Item(9, 3).Value = 0.001507075852714479
$ws.Cells.Item(9, 4).Value = 0.999764621257782
$ws.Cells.Item(10, 1).Value = 0.001244012615643442
$ws.Cells.Item(10, 2).Value = 0.9997853636741638
$ws.Cells.Item(10, 3).Value = 0.0007887334795668721
$ws.Cells.Item(10, 4).Value = 0.9996861219406128
$ws.Cells.Item(11, 1).Value = 0.001406561117619276
$ws.Cells.Item(11, 2).Value = 0.999666154384613
$ws.Cells.Item(11, 3).Value = 0.003313457826152444
$ws.Cells.Item(11, 4).Value = 0.9986660480499268
$ws.Cells.Item(12, 1).Value = 0.0008873771876096725
$ws.Cells.Item(12, 2).Value = 0.9998331069946289
$ws.Cells.Item(12, 3).Value = 0.003963550552725792
$ws.Cells.Item(12, 4).Value = 0.9981952309608459
$ws.Cells.Item(13, 1).Value = 0.001155850477516651
$ws.Cells.Item(13, 2).Value = 0.9997138381004333
$ws.Cells.Item(13, 3).Value = 0.001831115805543959
$ws.Cells.Item(13, 4).Value = 0.9996076822280884
$ws.Cells.Item(14, 1).Value = 0.001012247870676219
$ws.Cells.Item(14, 2).Value = 0.9997853636741638
$ws.Cells.Item(14, 3).Value = 0.0009772931225597858
$ws.Cells.Item(14, 4).Value = 0.9996861219406128
$ws.Cells.Item(15, 1).Value = 0.000583845132496208
$ws.Cells.Item(15, 2).Value = 0.9997376799583435
$ws.Cells.Item(15, 3).Value = 0.009900409728288651
$ws.Cells.Item(15, 4).Value = 0.9971751570701599
$ws.Cells.Item(16, 1).Value = 0.0006925089401192963
$ws.Cells.Item(16, 2).Value = 0.9998331069946289
$ws.Cells.Item(16, 3).Value = 0.005516128148883581
$ws.Cells.Item(16, 4).Value = 0.9976459741592407
$ws.Cells.Item(17, 1).Value = 0.0001008109538815916
$ws.Cells.Item(17, 2).Value = 0.9999761581420898
$ws.Cells.Item(17, 3).Value = 0.001042569754645228
$ws.Cells.Item(17, 4).Value = 0.9996076822280884
$ws.Cells.Item(18, 1).Value = 0.001048156060278416
$ws.Cells.Item(18, 2).Value = 0.9997138381004333
$ws.Cells.Item(18, 3).Value = 0.001986816059798002
$ws.Cells.Item(18, 4).Value = 0.9992938041687012
$ws.Cells.Item(19, 1).Value = 0.0006120206089690328
$ws.Cells.Item(19, 2).Value = 0.9999046325683594
$ws.Cells.Item(19, 3).Value = 0.000832011632155627
$ws.Cells.Item(19, 4).Value = 0.9999215602874756
$ws.Cells.Item(20, 1).Value = 0.00009467918425798416
$ws.Cells.Item(20, 2).Value = 0.9999761581420898
$ws.Cells.Item(20, 3).Value = 0.0004437766037881374
$ws.Cells.Item(20, 4).Value = 0.9999215602874756
$ws.Cells.Item(21, 1).Value = 0.0004156877985224128
$ws.Cells.Item(21, 2).Value = 0.9999046325683594
$ws.Cells.Item(21, 3).Value = 0.003781295381486416
$ws.Cells.Item(21, 4).Value = 0.9990583658218384
$ws.Cells.Item(22, 1).Value = 0.001289392355829477
$ws.Cells.Item(22, 2).Value = 0.999809205532074
$ws.Cells.Item(22, 3).Value = 0.0001973091129912063
$ws.Cells.Item(22, 4).Value = 0.9999215602874756
$ws.Cells.Item(23, 1).Value = 0.0005266459193080664
$ws.Cells.Item(23, 2).Value = 0.9998331069946289
$ws.Cells.Item(23, 3).Value = 0.001130411052145064
$ws.Cells.Item(23, 4).Value = 0.999764621257782
$ws.Cells.Item(24, 1).Value = 0.00007304515020223334
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = 0.0007031516288407147
$ws.Cells.Item(24, 4).Value = 0.999764621257782
$ws.Cells.Item(25, 1).Value = 0.0004401331534609199
$ws.Cells.Item(25, 2).Value = 0.9999046325683594
$ws.Cells.Item(25, 3).Value = 0.0006734869093634188
$ws.Cells.Item(25, 4).Value = 0.9999215602874756
$ws.Cells.Item(26, 1).Value = 0.000170375540619716
$ws.Cells.Item(26, 2).Value = 0.9999284744262695
$ws.Cells.Item(26, 3).Value = 0.00004496751353144646
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.0009175667073577642
$ws.Cells.Item(27, 2).Value = 0.999809205532074
$ws.Cells.Item(27, 3).Value = 0.00006813794607296586
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = 0.000694229151122272
$ws.Cells.Item(28, 2).Value = 0.9998569488525391
$ws.Cells.Item(28, 3).Value = 0.0008815351757220924
$ws.Cells.Item(28, 4).Value = 0.9994507431983948
$ws.Cells.Item(29, 1).Value = 0.0004813307605218142
$ws.Cells.Item(29, 2).Value = 0.9998807907104492
$ws.Cells.Item(29, 3).Value = 0.0001011676431517117
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = 0.0001832985290093347
$ws.Cells.Item(30, 2).Value = 0.9999761581420898
$ws.Cells.Item(30, 3).Value = 0.00007877962343627587
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(31, 1).Value = 0.0000309995302814059
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 0.0000402056029997766
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.00001307673846895341
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 0.00001268733194592642
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.0000345670196111314
$ws.Cells.Item(33, 2).Value = 0.9999761581420898
$ws.Cells.Item(33, 3).Value = 0.0003101641486864537
$ws.Cells.Item(33, 4).Value = 0.9999215602874756
$ws.Cells.Item(34, 1).Value = 0.0007797896978445351
$ws.Cells.Item(34, 2).Value = 0.9999284744262695
$ws.Cells.Item(34, 3).Value = 0.003792723640799522
$ws.Cells.Item(34, 4).Value = 0.9981167316436768
$ws.Cells.Item(35, 1).Value = 0.0006700505618937314
$ws.Cells.Item(35, 2).Value = 0.9998569488525391
$ws.Cells.Item(35, 3).Value = 0.0002147698687622324
$ws.Cells.Item(35, 4).Value = 0.9998430609703064
$ws.Cells.Item(36, 1).Value = 0.0004515814944170415
$ws.Cells.Item(36, 2).Value = 0.9999761581420898
$ws.Cells.Item(36, 3).Value = 0.00004451765926205553
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(37, 1).Value = 0.00007142062531784177
$ws.Cells.Item(37, 2).Value = 0.9999761581420898
$ws.Cells.Item(37, 3).Value = 0.000001488947418692987
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = 0.00002016609141719528
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0.00006911841774126515
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = 0.00004546483978629112
$ws.Cells.Item(39, 2).Value = 0.9999761581420898
$ws.Cells.Item(39, 3).Value = 0.00285354838706553
$ws.Cells.Item(39, 4).Value = 0.9994507431983948
$ws.Cells.Item(40, 1).Value = 0.001548729836940765
$ws.Cells.Item(40, 2).Value = 0.9998331069946289
$ws.Cells.Item(40, 3).Value = 0.00144458282738924
$ws.Cells.Item(40, 4).Value = 0.9996076822280884
$ws.Cells.Item(41, 1).Value = 0.0002306812384631485
$ws.Cells.Item(41, 2).Value = 0.9999046325683594
$ws.Cells.Item(41, 3).Value = 0.002107588341459632
$ws.Cells.Item(41, 4).Value = 0.9996076822280884
$ws.Cells.Item(42, 1).Value = 0.00004449577318155207
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 3).Value = 0.001186553854495287
$ws.Cells.Item(42, 4).Value = 0.9996076822280884
$ws.Cells.Item(43, 1).Value = 0.0009634911548346281
$ws.Cells.Item(43, 2).Value = 0.9998807907104492
$ws.Cells.Item(43, 3).Value = 0.00002136895272997208
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(44, 1).Value = 0.0001880240452010185
$ws.Cells.Item(44, 2).Value = 0.9999284744262695
$ws.Cells.Item(44, 3).Value = 0.0002064207656076178
$ws.Cells.Item(44, 4).Value = 0.9999215602874756
$ws.Cells.Item(45, 1).Value = 0.0006568202516064048
$ws.Cells.Item(45, 2).Value = 0.9999523162841797
$ws.Cells.Item(45, 3).Value = 0.0005454533384181559
$ws.Cells.Item(45, 4).Value = 0.9998430609703064
$ws.Cells.Item(46, 1).Value = 0.0009652992594055831
$ws.Cells.Item(46, 2).Value = 0.999809205532074
$ws.Cells.Item(46, 3).Value = 0.0007611916516907513
$ws.Cells.Item(46, 4).Value = 0.9996861219406128
$ws.Cells.Item(47, 1).Value = 0.00006486885104095563
$ws.Cells.Item(47, 2).Value = 0.9999761581420898
$ws.Cells.Item(47, 3).Value = 0.0008561391150578856
$ws.Cells.Item(47, 4).Value = 0.9996861219406128
$ws.Cells.Item(48, 1).Value = 0.0001817394804675132
$ws.Cells.Item(48, 2).Value = 0.9999284744262695
$ws.Cells.Item(48, 3).Value = 0.0004921953659504652
$ws.Cells.Item(48, 4).Value = 0.9999215602874756
$ws.Cells.Item(49, 1).Value = 0.0004948212299495935
$ws.Cells.Item(49, 2).Value = 0.9999761581420898
$ws.Cells.Item(49, 3).Value = 0.02639380097389221
$ws.Cells.Item(49, 4).Value = 0.9890144467353821
$ws.Cells.Item(50, 1).Value = 0.0001238637632923201
$ws.Cells.Item(50, 2).Value = 0.9999761581420898
$ws.Cells.Item(50, 3).Value = 0.00007642973650945351
$ws.Cells.Item(50, 4).Value = 0.9999215602874756
$ws.Cells.Item(51, 1).Value = 0.000006831506652815733
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = 0.0001523956598248333
$ws.Cells.Item(51, 4).Value = 0.9999215602874756

$wb.Save()
